$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values per the diff
$ws.Range("B28").Value = 381
$ws.Range("B38").Value = 537

# Add new row 39 for week 38
$ws.Range("A39").Value = 38
$ws.Range("B39").Value = 612
